$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.856.30"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").Value = "1.906.71"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").Value = "313.65"
$ws.Range("E5").Value = "  -1.00%  "

# Row 6
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("D7").Value = "0.5026"
$ws.Range("E7").Value = "  +4.13%  "

# Row 8
$ws.Range("D8").Value = "0.3817"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "0.07282"
$ws.Range("E9").Value = "  -1.10%  "

# Row 10
$ws.Range("D10").Value = "0.9084"
$ws.Range("E10").Value = "  -2.88%  "

# Row 11
$ws.Range("D11").Value = "20.87"
$ws.Range("E11").Value = "  +0.17%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07671"
$ws.Range("E12").Value = "  -1.65%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.906.78"
$ws.Range("E13").Value = "  -0.34%  "

# Row 14
$ws.Range("D14").Value = "5.478"
$ws.Range("E14").Value = "  -0.47%  "

# Row 15
$ws.Range("D15").Value = "91.86"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.26%  "

# Row 17
$ws.Range("D17").Value = "0.000008714"
$ws.Range("E17").Value = "  -1.44%  "

# Row 18
$ws.Range("E18").Value = "  -0.28%  "

# Row 19
$ws.Range("D19").Value = "27.895.48"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20
$ws.Range("E20").Value = "  -1.60%  "

# Row 21
$ws.Range("D21").Value = "5.168"
$ws.Range("E21").Value = "  -0.24%  "

# Row 22
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  -0.95%  "

# Row 23
$ws.Range("D23").Value = "6.583"
$ws.Range("E23").Value = "  -0.89%  "

# Row 24
$ws.Range("D24").Value = "153.90"
$ws.Range("E24").Value = "  -1.14%  "

# Row 25
$ws.Range("D25").Value = "1.876"
$ws.Range("E25").Value = "  -2.53%  "

# Row 26
$ws.Range("D26").Value = "2.218"
$ws.Range("E26").Value = "  +4.58%  "

# Row 27
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  -0.82%  "

# Row 28
$ws.Range("D28").Value = "115.25"
$ws.Range("E28").Value = "  -1.19%  "

# Row 29
$ws.Range("D29").Value = "4.906"
$ws.Range("E29").Value = "  -1.23%  "

# Row 30
$ws.Range("D30").Value = "0.09008"
$ws.Range("E30").Value = "  +0.51%  "

# Row 31
$ws.Range("D31").Value = "3.214"
$ws.Range("E31").Value = "  -2.94%  "

# Row 32
$ws.Range("D32").Value = "1.226"
$ws.Range("E32").Value = "  -1.97%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.652"
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7614"
$ws.Range("E34").Value = "  -1.72%  "

# Row 35
$ws.Range("D35").Value = "0.02061"
$ws.Range("E35").Value = "  +0.29%  "

# Row 36
$ws.Range("D36").Value = "2.531"
$ws.Range("E36").Value = "  -4.67%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.092"
$ws.Range("E37").Value = "  -1.78%  "

# Row 38
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.5566"
$ws.Range("E38").Value = "  +1.66%  "

# Row 39
$ws.Range("D39").Value = "3.023"
$ws.Range("E39").Value = "  +0.90%  "

# Row 40
$ws.Range("D40").Value = "0.05250"
$ws.Range("E40").Value = "  -1.10%  "

# Row 41
$ws.Range("D41").Value = "6.907"
$ws.Range("E41").Value = "  -1.78%  "

# Row 42
$ws.Range("D42").Value = "8.483"
$ws.Range("E42").Value = "  -0.35%  "

# Row 43
$ws.Range("D43").Value = "0.1511"
$ws.Range("E43").Value = "  -1.11%  "

# Row 44
$ws.Range("D44").Value = "111.12"
$ws.Range("E44").Value = "  +2.39%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "10.62"
$ws.Range("E45").Value = "  -1.28%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.4825"
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("E47").Value = "  -0.40%  "

# Row 48
$ws.Range("D48").Value = "1.626"
$ws.Range("E48").Value = "  -1.64%  "

# Row 49
$ws.Range("D49").Value = "67.38"
$ws.Range("E49").Value = "  -1.17%  "

# Row 50
$ws.Range("D50").Value = "0.06069"
$ws.Range("E50").Value = "  -0.29%  "

# Row 51
$ws.Range("D51").Value = "0.9029"
$ws.Range("E51").Value = "  +0.23%  "
